$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.169.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.39%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.556.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.03%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.33%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.555.43'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.91%  '

# Row 10
$ws.Range("E10").Value = '  +5.68%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.21%  '

# Row 12
$ws.Range("E12").Value = '  +3.46%  '

# Row 13
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.76%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.95%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.158.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.00%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.550.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.85%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.530.64'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.01%  '

# Row 18
$ws.Range("E18").Value = '  +0.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.07%  '

# Row 20
$ws.Range("E20").Value = '  +5.85%  '

# Row 21
$ws.Range("E21").Value = '  +10.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '456.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.640'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.99%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.50%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.00%  '

# Row 26
$ws.Range("E26").Value = '  +1.85%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.694.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.88%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.08'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.33%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.99%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.172'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.47%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.12%  '

# Row 34
$ws.Range("E34").Value = '  +4.47%  '

# Row 35
$ws.Range("E35").Value = '  +1.47%  '

# Row 36
$ws.Range("E36").Value = '  +3.66%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.551.17'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.08%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.25'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.39%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.56%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '178.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.51%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0918'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.45%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.15%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.29%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +15.73%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.895'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.41%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.67%  '

# Row 48
$ws.Range("E48").Value = '  +7.25%  '

# Row 49
$ws.Range("E49").Value = '  +3.89%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.27%  '

# Row 51
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.261'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.72%  '
